$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B85 was stored as text "1"; convert it to a real number (value 1)
$ws.Range("B85").Value = 1

# Append new row 86 with annotation data for Ying Tang
$ws.Range("A86").Value = "Ying Tang"

# B86 keeps the numeric-looking value "3" stored as text, like the source data
$ws.Range("B86").NumberFormat = "@"
$ws.Range("B86").Value = "3"
$ws.Range("B86").ClearFormats()

$ws.Range("C86").Value = "无"
$ws.Range("D86").Value = "SMY"
$ws.Range("E86").Value = "RES"
$ws.Range("F86").Value = "94664fc5-740b-497e-9f27-9fbb6b5fbbdd"
$ws.Range("G86").Value = "TT0bFo9VZpFWg_annotated.xlsx"
$ws.Range("H86").Value = "The net gets bigger, yet keeps underfitting the training set."
